$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.146142244338989
$ws.Range("B1").Value = 2.149334907531738
$ws.Range("C1").Value = 2.91516375541687
$ws.Range("D1").Value = 3.486250638961792
$ws.Range("E1").Value = 1.637382626533508
